# Apply the "Fruta / hortaliza, semanal" update.
# The data rows got re-ordered (weekly date re-shuffle), which in terms of
# the underlying cell values amounts to swapping the "Fecha" (D), "Volumen"
# (M), "Precio minimo" (N), "Precio maximo" (O), "Precio promedio
# ponderado" (P) and "Precio $/Kg" (S) values between certain row pairs,
# while every other column (A, B, C, E, F, G, H, I, J, K, L, Q, R, T) stays
# identical across all rows and therefore needs no change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowValues {
    param([int]$RowA, [int]$RowB)

    $columns = @("D", "M", "N", "O", "P", "S")

    foreach ($col in $columns) {
        $cellA = $ws.Range("$col$RowA")
        $cellB = $ws.Range("$col$RowB")

        $valueA = $cellA.Value2
        $valueB = $cellB.Value2

        $cellA.Value2 = $valueB
        $cellB.Value2 = $valueA
    }
}

# Pairwise swaps derived from the diff (row 4 is untouched).
Swap-RowValues 3 8
Swap-RowValues 5 7
Swap-RowValues 6 10
Swap-RowValues 9 11
